$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115 — this shifts the existing rows 115..179
# down to 116..180 (Excel's normal "insert row" semantics), matching the
# diff where every record previously at row N (115<=N<=179) now lives at
# row N+1, and a brand new record is written into the freed row 115.
$ws.Rows("115:115").Insert()

# Populate the newly inserted row 115 with the new record.
$ws.Cells.Item(115, 1).Value  = 7
$ws.Cells.Item(115, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value  = "Ñuble"
$ws.Cells.Item(115, 4).Value  = 44438
$ws.Cells.Item(115, 5).Value  = 16
$ws.Cells.Item(115, 6).Value  = 100114001
$ws.Cells.Item(115, 7).Value  = "Papa"
$ws.Cells.Item(115, 8).Value  = "Patagonia"
$ws.Cells.Item(115, 9).Value  = "1a (guarda)"
$ws.Cells.Item(115, 10).Value = 300
$ws.Cells.Item(115, 11).Value = 6500
$ws.Cells.Item(115, 12).Value = 7000
$ws.Cells.Item(115, 13).Value = 6750
$ws.Cells.Item(115, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(115, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(115, 16).Value = 270
$ws.Cells.Item(115, 17).Value = 25
$ws.Cells.Item(115, 18).Value = "Hortaliza"
